$d = $word.ActiveDocument

# --- Paragraph 1: "Phonegap plugin - Quick Look" -> bump font size to 24pt (sz/szCs = 48) ---
$p1 = $d.Paragraphs(1)
$p1.Range.Font.Size = 24
$p1.Range.Font.SizeBi = 24

# --- Paragraph 2: empty paragraph -> still needs the paragraph-mark run formatting (pPr/rPr) ---
# Setting Font.* directly on an empty paragraph's range fails in this runtime ("collection
# member does not exist") because there is no run to carry the formatting.  Work around it by
# temporarily inserting a character, formatting the (now non-empty) paragraph, and then
# deleting just that character again - the pPr/rPr mark formatting survives the deletion.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertAfter("X")
$p2b = $d.Paragraphs(2)
$p2b.Range.Font.Size = 24
$p2b.Range.Font.SizeBi = 24
$p2c = $d.Paragraphs(2)
$tempChar = $d.Range($p2c.Range.Start, $p2c.Range.Start + 1)
$tempChar.Delete()

# --- Paragraph 3: "Demo of PDF Document" -> "Demo of Word Document", bigger font ---
# First mark the "PDF" run as distinct (Bold toggle) *before* rewriting its text so the
# engine does not silently coalesce it into the neighbouring "Demo of " run when the text
# is replaced (adjacent runs with identical formatting get merged on a plain text write).
$p3 = $d.Paragraphs(3)
$pdfStart = $p3.Range.Start + 8
$pdfEnd = $pdfStart + 3
$pdfRange = $d.Range($pdfStart, $pdfEnd)
$pdfRange.Font.Bold = 1

$pdfRange2 = $d.Range($pdfStart, $pdfEnd)
$pdfRange2.Text = "Word"

# Now size the whole paragraph (adds pPr/rPr mark formatting + sizes the two untouched runs).
$p3b = $d.Paragraphs(3)
$p3b.Range.Font.Size = 24
$p3b.Range.Font.SizeBi = 24

# Finally restore the "Word" run to plain (un-bold) and give it the same bigger size.
$wordStart = $pdfStart
$wordEnd = $wordStart + 4
$wordRange = $d.Range($wordStart, $wordEnd)
$wordRange.Font.Bold = 0
$wordRange.Font.Size = 24
$wordRange.Font.SizeBi = 24
